# Update countries & provincias Spain
# - Nepal moves up in ranking (new data inserted at row 122), pushing
#   Principado de Andorra .. Republica del Chad down by one row.
# - Madagascar moves up one rank, swapping with Jamaica.
# - Several countries' case counts refreshed (India, Catar, Finlandia, Vietnam,
#   and the whole Andorra..Chad / Jamaica..Madagascar block).
# - "Datos actualizados" timestamp bumped from 13:05 to 13:35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp row ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 13:35"

# --- India (row 13) ---
$ws.Range("B13").Value = 146376
$ws.Range("C13").Value = 1426
$ws.Range("E13").Value = 81038

# --- Catar (row 23) ---
$ws.Range("B23").Value = 47207
$ws.Range("C23").Value = 1742
$ws.Range("D23").Value = 11844
$ws.Range("E23").Value = 35335
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 28

# --- Finlandia (row 67) ---
$ws.Range("E67").Value = 1216
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 312

# --- Nepal jumps into the ranking at row 122, shifting Andorra..Chad down ---
$ws.Range("A122").Value = "Nepal"
$ws.Range("B122").Value = 772
$ws.Range("C122").Value = 90
$ws.Range("D122").Value = 155
$ws.Range("E122").Value = 613
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 4

$ws.Range("A123").Value = "Principado de Andorra"
$ws.Range("B123").Value = 763
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 663
$ws.Range("E123").Value = 49
$ws.Range("H123").Value = 51

$ws.Range("A124").Value = "Sierra Leona"
$ws.Range("B124").Value = 735
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 293
$ws.Range("E124").Value = 400
$ws.Range("H124").Value = 42

$ws.Range("A125").Value = "Georgia"
$ws.Range("B125").Value = 732
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = 537
$ws.Range("E125").Value = 183
$ws.Range("H125").Value = 12

$ws.Range("A126").Value = "Crucero"
$ws.Range("B126").Value = 712
$ws.Range("D126").Value = 651
$ws.Range("E126").Value = 48
$ws.Range("H126").Value = 13

$ws.Range("A127").Value = "Jordania"
$ws.Range("B127").Value = 711
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 479
$ws.Range("E127").Value = 223
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 9

$ws.Range("A128").Value = "Etiopia"
$ws.Range("B128").Value = 701
$ws.Range("C128").Value = 46
$ws.Range("D128").Value = 167
$ws.Range("E128").Value = 528
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 6

$ws.Range("A129").Value = "Republica del Chad"
$ws.Range("B129").Value = 687
$ws.Range("D129").Value = 244
$ws.Range("E129").Value = 382
$ws.Range("H129").Value = 61

# --- Madagascar overtakes Jamaica at row 133 ---
$ws.Range("A133").Value = "Madagascar"
$ws.Range("B133").Value = 586
$ws.Range("C133").Value = 44
$ws.Range("D133").Value = 147
$ws.Range("E133").Value = 437
$ws.Range("H133").Value = 2

$ws.Range("A134").Value = "Jamaica"
$ws.Range("B134").Value = 556
$ws.Range("C134").Value = 4
$ws.Range("D134").Value = 238
$ws.Range("E134").Value = 309
$ws.Range("H134").Value = 9

# --- Vietnam (row 146) ---
$ws.Range("B146").Value = 327
$ws.Range("C146").Value = 1
$ws.Range("E146").Value = 55
